$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.587.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.647.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.534"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.89%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.63"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.653.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("E14").Value = "  +3.77%  "
$ws.Range("E15").Value = "  -2.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.545.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.58%  "
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  -3.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.55%  "
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.74%  "
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0485"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("E33").Value = "  +2.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.424.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.569"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  -3.70%  "
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.816"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.789.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0995"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.34%  "
